# The commit sets the whole document's font size to 18pt (w:sz/w:szCs =
# 36 half-points), including the paragraph-mark run properties (w:pPr/
# w:rPr) for every paragraph - even the empty one.
$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range

    if ($r.Characters.Count -le 1) {
        # Paragraph has no run (just the paragraph mark) - Font.Size on
        # such a range is a no-op because there is no run to carry the
        # formatting. Temporarily insert a character so the size can be
        # stamped (which also stamps the paragraph mark's rPr), then
        # remove the character again, leaving the paragraph-mark
        # formatting behind.
        $r.InsertBefore("X")
        $r.Font.Size = 18
        $r.Font.SizeBi = 18
        $d.Range($r.Start, $r.Start + 1).Text = ""
    } else {
        $r.Font.Size = 18
        $r.Font.SizeBi = 18
    }
}
